# Apply the data refresh recorded in the commit "Update gh-pages to output
# generated at 456a3b4" to 杭州-漫展信息.xlsx.
#
# The workbook has 4 worksheets: 展览 (exhibitions), 演出 (performances),
# 本地生活 (local life) and 全部类型 (all types, a merged/aggregated view).
# Most of the edits are simple refreshes of the "想去人数" (want-to-go count,
# column F) numbers. In addition:
#   - 展览!G12 flips from a price (65) to the text "不可售" (not for sale)
#   - 全部类型 rows 14-21 get new content, shifted because the "初始之音响彻
#     未来" show became unavailable and dropped out of the aggregated sheet,
#     while the "COMIC GALAXY" show newly appears at the end of that block.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------
# Sheet 展览 - simple numeric refreshes in column F
# ---------------------------------------------------------------------
$ws1.Range("F3").Value  = 409
$ws1.Range("F4").Value  = 1154
$ws1.Range("F5").Value  = 42
$ws1.Range("F7").Value  = 28
$ws1.Range("F8").Value  = 1070
$ws1.Range("F10").Value = 350
$ws1.Range("F11").Value = 423
$ws1.Range("F17").Value = 495
$ws1.Range("F18").Value = 635
$ws1.Range("F19").Value = 5647
$ws1.Range("F21").Value = 1571
$ws1.Range("F22").Value = 372
$ws1.Range("F23").Value = 31
$ws1.Range("F24").Value = 21
$ws1.Range("F25").Value = 4979
$ws1.Range("F26").Value = 121
$ws1.Range("F28").Value = 1508
$ws1.Range("F31").Value = 654
$ws1.Range("F32").Value = 75

# G12: ticket price becomes unavailable for purchase
$ws1.Range("G12").Value = "不可售"

# ---------------------------------------------------------------------
# Sheet 演出 - simple numeric refreshes in column F
# ---------------------------------------------------------------------
$ws2.Range("F4").Value  = 13
$ws2.Range("F5").Value  = 144
$ws2.Range("F8").Value  = 117
$ws2.Range("F12").Value = 1
$ws2.Range("F13").Value = 15

# ---------------------------------------------------------------------
# Sheet 本地生活 - simple numeric refreshes in column F
# ---------------------------------------------------------------------
$ws3.Range("F2").Value = 9392
$ws3.Range("F4").Value = 2136

# ---------------------------------------------------------------------
# Sheet 全部类型 - simple numeric refreshes in column F
# (these rows correspond unchanged to rows already updated above)
# ---------------------------------------------------------------------
$ws4.Range("F2").Value  = 9392
$ws4.Range("F4").Value  = 2136
$ws4.Range("F6").Value  = 409
$ws4.Range("F7").Value  = 1154
$ws4.Range("F8").Value  = 42
$ws4.Range("F10").Value = 28
$ws4.Range("F11").Value = 1070
$ws4.Range("F12").Value = 350
$ws4.Range("F13").Value = 423
$ws4.Range("F22").Value = 640
$ws4.Range("F23").Value = 5647
$ws4.Range("F25").Value = 1571
$ws4.Range("F28").Value = 372
$ws4.Range("F30").Value = 1
$ws4.Range("F31").Value = 4982
$ws4.Range("F32").Value = 121
$ws4.Range("F34").Value = 1508
$ws4.Range("F37").Value = 654
$ws4.Range("F38").Value = 75
$ws4.Range("F39").Value = 15

# ---------------------------------------------------------------------
# Sheet 全部类型 - rows 14-21 get fully rewritten content: the
# "初始之音响彻未来" event (old row 14) drops out because it became
# unavailable, rows 15-20 shift up into rows 14-19... in effect every
# row from 14 to 20 takes on the content that used to belong to the
# next row, and row 21 is populated with the brand new "COMIC GALAXY"
# entry.
# ---------------------------------------------------------------------

# Row 14: 杭州·初音未来17周年生日派对 & 音链视窗同人共创only
$ws4.Range("C14").Value = "杭州·初音未来17周年生日派对 & 音链视窗同人共创only"
$ws4.Range("D14").Value = "金惠路1128号西区 杭州金迪大酒店"
$ws4.Range("E14").Value = "2024.08.31 12:00-08.31 20:00"
$ws4.Range("F14").Value = 312
$ws4.Range("G14").Value = 39
$ws4.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=90372"
$ws4.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202408/j61I0I7n1722925055083.jpeg"

# Row 15: 杭州·音乐番+only
$ws4.Range("C15").Value = "杭州·音乐番+only"
$ws4.Range("D15").Value = "康候圣街99号 顺丰创新中心"
$ws4.Range("E15").Value = "2024.08.31 10:00-08.31 18:00"
$ws4.Range("F15").Value = 356
$ws4.Range("G15").Value = 78
$ws4.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=88899"
$ws4.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202408/lxW52TpT1724228135568.jpeg"

# Row 16: 杭州·次元萌友会【免费展会】
$ws4.Range("B16").Value = "2024-09-06"
$ws4.Range("C16").Value = "杭州·次元萌友会【免费展会】"
$ws4.Range("D16").Value = "祥泰街398号 杭州万融城"
$ws4.Range("E16").Value = "2024.09.06 10:00-09.08 21:00"
$ws4.Range("F16").Value = 31
$ws4.Range("G16").Value = 20
$ws4.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=90896"
$ws4.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202408/snpy8ATR1723793956830.png"

# Row 17: 杭州·DNP01综合同人展X【昼夜星逐】泛VOCALOID专场
$ws4.Range("B17").Value = "2024-09-07"
$ws4.Range("C17").Value = "杭州·DNP01综合同人展X【昼夜星逐】泛VOCALOID专场"
$ws4.Range("D17").Value = "观澜路钱江世纪公园d区1幢 杭州世纪雷迪森庄园酒店"
$ws4.Range("E17").Value = "2024.09.07 12:00-09.07 20:00"
$ws4.Range("F17").Value = 64
$ws4.Range("G17").Value = 20
$ws4.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=90587"
$ws4.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202408/x9rSjkDf1723208824749.jpeg"

# Row 18: 杭州·《卡农》永恒经典名曲音乐会
$ws4.Range("C18").Value = "杭州·《卡农》永恒经典名曲音乐会"
$ws4.Range("D18").Value = "曙光路31号 浙江音乐厅"
$ws4.Range("E18").Value = "2024.09.07 19:30-09.07 21:00"
$ws4.Range("F18").Value = 2
$ws4.Range("G18").Value = 100
$ws4.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=85894"
$ws4.Range("I18").Value = "//i2.hdslb.com/bfs/openplatform/202405/3jz9YpaW1716100738530.jpeg"

# Row 19: 杭州·红楼梦·主题演绎国风音乐会《梦寻红楼》
$ws4.Range("C19").Value = "杭州·红楼梦·主题演绎国风音乐会《梦寻红楼》"
$ws4.Range("D19").Value = "望梅路与汀兰路交叉口向南100米 杭州临平大剧院（原余杭大剧院）"
$ws4.Range("E19").Value = "2024.09.07 15:00-09.07 16:30"
$ws4.Range("F19").Value = 17
$ws4.Range("G19").Value = 100
$ws4.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=89257"
$ws4.Range("I19").Value = "//i2.hdslb.com/bfs/openplatform/202407/tkm6AHo71720572975141.jpeg"

# Row 20: 杭州·2024吉卜力动漫音乐原版歌手交响音乐会
$ws4.Range("B20").Value = "2024-09-15"
$ws4.Range("C20").Value = "杭州·2024吉卜力动漫音乐原版歌手交响音乐会"
$ws4.Range("D20").Value = "魔方剧院三楼 武林广场"
$ws4.Range("E20").Value = "2024.09.15 19:30-09.15 21:00"
$ws4.Range("F20").Value = 4
$ws4.Range("G20").Value = 280
$ws4.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=89692"
$ws4.Range("I20").Value = "//i2.hdslb.com/bfs/openplatform/202407/jzRpaVmw1721703321314.png"

# Row 21: 杭州·2024首届COMIC GALAXY次元盛典 (newly appended)
$ws4.Range("C21").Value = "杭州·2024首届COMIC GALAXY次元盛典"
$ws4.Range("D21").Value = "长江南路336号 白马湖国际会展中心"
$ws4.Range("E21").Value = "2024.09.15 09:30-09.17 17:30"
$ws4.Range("F21").Value = 495
$ws4.Range("G21").Value = 88
$ws4.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=90433"
$ws4.Range("I21").Value = "//i0.hdslb.com/bfs/openplatform/202408/teoBMbzd1723019674766.png"
